# Fixed Signal 11,12,13,14,15 unit
#
# The "Unit" column (F) for the Radar Sensor signals with Signal ID
# 11-15 (rows 15-19 on the "CommunicationMatrix" sheet) was inconsistent
# ("km/h", "m", "Pixel"). All five are corrected to the single, lower
# case unit "pixel".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CommunicationMatrix")

# Signal 11 - Longitudinal RVX (row 15): was "km/h"
$ws.Range("F15").Value = "pixel"

# Signal 12 - Longitudinal EGO (row 16): was "m"
$ws.Rows.Item(16).RowHeight = 12.8
$ws.Range("F16").Value = "pixel"

# Signal 13 - Lateral RVY (row 17): was "km/h"
$ws.Rows.Item(17).RowHeight = 12.8
$ws.Range("F17").Value = "pixel"

# Signal 14 - Lateral EGO (row 18): was "m"
$ws.Rows.Item(18).RowHeight = 12.8
$ws.Range("F18").Value = "pixel"

# Signal 15 - Object Size (row 19): was "Pixel" (capitalised)
$ws.Range("F19").Value = "pixel"

# Make "CommunicationMatrix" the active sheet/tab with F15 selected,
# matching the author's final cursor position (activeTab 0 -> 1,
# tabSelected moves from "Version & History" to "CommunicationMatrix").
$ws.Activate()
$ws.Range("F15").Select()
